# Update crypto price/volume data as of the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column holds price text that must stay literal text (it uses a
# dotted thousands format like "58.007.27" and fixed decimals like
# "1.00"), so force the Text format before writing the new values.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '58.007.27'
$ws.Range('E2').Value = '  -4.09%  '

$ws.Range('D3').Value = '2.971.83'
$ws.Range('E3').Value = '  -1.13%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').Value = '557.94'
$ws.Range('E5').Value = '  -3.91%  '

$ws.Range('D6').Value = '133.50'
$ws.Range('E6').Value = '  +5.25%  '

$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('D8').Value = '0.514'
$ws.Range('E8').Value = '  +3.02%  '

$ws.Range('D9').Value = '2.966.80'
$ws.Range('E9').Value = '  -1.22%  '

$ws.Range('E10').Value = '  -3.06%  '

$ws.Range('D11').Value = '4.88'
$ws.Range('E11').Value = '  -5.28%  '

$ws.Range('D12').Value = '0.449'
$ws.Range('E12').Value = '  +1.79%  '

$ws.Range('D13').Value = '0.0000224'
$ws.Range('E13').Value = '  -0.31%  '

$ws.Range('D14').Value = '33.03'
$ws.Range('E14').Value = '  +1.27%  '

$ws.Range('E15').Value = '  +0.61%  '

$ws.Range('D16').Value = '3.462.41'
$ws.Range('E16').Value = '  -1.02%  '

$ws.Range('D17').Value = '6.89'
$ws.Range('E17').Value = '  +7.82%  '

$ws.Range('D18').Value = '2.970.82'
$ws.Range('E18').Value = '  -1.17%  '

$ws.Range('D19').Value = '57.920.99'
$ws.Range('E19').Value = '  -3.98%  '

$ws.Range('D20').Value = '420.13'
$ws.Range('E20').Value = '  -2.81%  '

$ws.Range('D21').Value = '13.23'
$ws.Range('E21').Value = '  +0.43%  '

$ws.Range('D22').Value = '0.688'
$ws.Range('E22').Value = '  +3.57%  '

$ws.Range('D23').Value = '7.00'
$ws.Range('E23').Value = '  -0.69%  '

$ws.Range('D24').Value = '13.12'
$ws.Range('E24').Value = '  +1.81%  '

$ws.Range('D25').Value = '79.66'
$ws.Range('E25').Value = '  +0.23%  '

$ws.Range('E26').Value = '  -0.14%  '

$ws.Range('E27').Value = '  +0.13%  '

$ws.Range('D28').Value = '2.50'
$ws.Range('E28').Value = '  -2.44%  '

$ws.Range('D29').Value = '7.59'
$ws.Range('E29').Value = '  +3.53%  '

$ws.Range('E30').Value = '  +5.52%  '

$ws.Range('D31').Value = '25.31'
$ws.Range('E31').Value = '  -0.31%  '

$ws.Range('D32').Value = '6.09'
$ws.Range('E32').Value = '  -1.46%  '

$ws.Range('E33').Value = '  +6.32%  '

$ws.Range('D34').Value = '2.14'
$ws.Range('E34').Value = '  -1.03%  '

$ws.Range('D35').Value = '5.66'
$ws.Range('E35').Value = '  +0.77%  '

$ws.Range('D36').Value = '0.942'
$ws.Range('E36').Value = '  -1.26%  '

$ws.Range('D37').Value = '0.0₃0696'
$ws.Range('E37').Value = '  +4.47%  '

$ws.Range('D38').Value = '48.72'
$ws.Range('E38').Value = '  -2.66%  '

$ws.Range('D39').Value = '8.59'
$ws.Range('E39').Value = '  +2.20%  '

$ws.Range('D40').Value = '2.58'
$ws.Range('E40').Value = '  +3.79%  '

# Rows 41/42 swap rank order: VeChain moves above Bittensor with refreshed
# price/volume figures.
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.0351'
$ws.Range('E41').Value = '  -2.64%  '

$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = '381.23'
$ws.Range('E42').Value = '  -0.90%  '

$ws.Range('E43').Value = '  -1.80%  '

$ws.Range('D44').Value = '2.687.72'
$ws.Range('E44').Value = '  +1.16%  '

$ws.Range('E45').Value = '  -0.03%  '

$ws.Range('D46').Value = '0.243'
$ws.Range('E46').Value = '  +2.70%  '

$ws.Range('D47').Value = '122.12'
$ws.Range('E47').Value = '  +3.25%  '

$ws.Range('E48').Value = '  +2.59%  '

$ws.Range('D49').Value = '1.99'
$ws.Range('E49').Value = '  -1.43%  '

$ws.Range('D50').Value = '23.58'
$ws.Range('E50').Value = '  -0.87%  '

$ws.Range('D51').Value = '2.02'
$ws.Range('E51').Value = '  -0.49%  '
